$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Stages:" legend list (column J) updates ---
$ws.Range("J7").Value = "Blocking"
$ws.Range("J8").Value = "Colour/Silhouette"
$ws.Range("J9").Value = "Definition"
$ws.Range("J10").Value = "Light"
$ws.Range("J11").Value = "Detail"
$ws.Range("J12").Value = "REDESIGN/RESIZE"
$ws.Range("J13").Value = "Ready for Placeholder"
$ws.Range("J14").Value = "Complete"

# --- Olive tree (row 10) update ---
$ws.Range("E8").Value = "Resize"
$ws.Range("E10").Value = "Ready for Placeholder"

# --- New Turkish Pine sprite (row 11 - phoenix theophrasti / Cretan Date Palm row) ---
$ws.Range("D11").Value = "Jayden"
$ws.Range("E11").Value = "Blocking"
$ws.Range("F11").Value = 1

# --- Row 12 (pinus brutia / Turkish Pine) new sprite entry ---
$ws.Range("D12").Value = "Jayden"
$ws.Range("E12").Value = "Light"

# --- N/A fills for Level Background section (rows 27-30, column C) ---
$ws.Range("C27").Value = "N/A"
$ws.Range("C27").Style = $ws.Range("C26").Style
$ws.Range("C28").Value = "N/A"
$ws.Range("C28").Style = $ws.Range("C26").Style
$ws.Range("C29").Value = "N/A"
$ws.Range("C29").Style = $ws.Range("C26").Style
$ws.Range("C30").Value = "N/A"
$ws.Range("C30").Style = $ws.Range("C26").Style

# --- Sarissa renamed to Sarissa (Spear) ---
$ws.Range("C49").Value = "Sarissa (Spear)"

# --- Column E width ---
# Target stored width is 43.28515625 character-units; the runtime quantizes
# ColumnWidth assignments to 1/6 increments, so 42.5 is the closest input
# that reproduces the nearest achievable stored width (43.33333...).
$ws.Range("E1").ColumnWidth = 42.5

# --- Selection ---
$ws.Range("E19").Select()
